$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume(1h) changes, and some row
# reassignments where a coin dropped out of the top list and a new one
# entered at the bottom).
# Column D ("Price") cells are stored as plain text in this sheet (e.g.
# "1.00", "27.30", "18.40"), so each D write forces text format first to
# stop Excel from auto-coercing the value into a number and dropping
# trailing zeros, then restores the "Normal" style so no stray number
# format is left applied to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.920.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.583.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.045.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.775.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.581.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.708.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.94%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0821"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "467.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.81%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "175.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.400"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.00%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "157.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.27%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.634"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0541"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0965"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
